# Apply cryptos list update per diff (cells B/C/D/E for affected rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.472.97"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "3.335.15"
$ws.Range("E3").Value = "  -0.17%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'190.10"
$ws.Range("E5").Value = "  +4.24%  "

$ws.Range("D6").Value = "'566.07"
$ws.Range("E6").Value = "  +1.07%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("E8").Value = "  -0.49%  "

$ws.Range("D9").Value = "3.329.17"
$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("E10").Value = "  +0.73%  "

$ws.Range("D12").Value = "'47.96"
$ws.Range("E12").Value = "  +0.87%  "

$ws.Range("E13").Value = "  +2.93%  "

$ws.Range("D14").Value = "'8.71"
$ws.Range("E14").Value = "  +0.69%  "

$ws.Range("D15").Value = "3.865.90"
$ws.Range("E15").Value = "  -0.30%  "

$ws.Range("D16").Value = "'606.99"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "66.452.10"
$ws.Range("E17").Value = "  +0.22%  "

$ws.Range("D18").Value = "'18.12"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("E19").Value = "  +1.15%  "

$ws.Range("D20").Value = "3.335.67"
$ws.Range("E20").Value = "  -0.32%  "

$ws.Range("D21").Value = "'11.18"
$ws.Range("E21").Value = "  -2.48%  "

$ws.Range("D22").Value = "'0.916"
$ws.Range("E22").Value = "  +0.99%  "

$ws.Range("D23").Value = "'18.99"
$ws.Range("E23").Value = "  +12.33%  "

$ws.Range("D24").Value = "'5.21"
$ws.Range("E24").Value = "  +3.16%  "

$ws.Range("D25").Value = "'101.51"
$ws.Range("E25").Value = "  +1.56%  "

$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'5.99"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").Value = "'2.77"
$ws.Range("E28").Value = "  +3.13%  "

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'9.77"
$ws.Range("E29").Value = "  +4.64%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").Value = "'8.74"
$ws.Range("E30").Value = "  +0.18%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").Value = "'30.61"
$ws.Range("E31").Value = "  +0.80%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'6.81"
$ws.Range("E32").Value = "  +8.72%  "

$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").Value = "'4.05"
$ws.Range("E33").Value = "  +6.92%  "

$ws.Range("D34").Value = "'565.64"
$ws.Range("E34").Value = "  +3.19%  "

$ws.Range("B35").Value = "Cosmos"
$ws.Range("C35").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D35").Value = "'11.16"
$ws.Range("E35").Value = "  +1.02%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = "  +1.63%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").Value = "'57.36"
$ws.Range("E37").Value = "  -0.25%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "3.718.10"
$ws.Range("E38").Value = "  -3.87%  "

$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0733"
$ws.Range("E40").Value = "  +2.25%  "

$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'34.17"
$ws.Range("E41").Value = "  +6.26%  "

$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.131"
$ws.Range("E42").Value = "  +4.97%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'3.32"
$ws.Range("E43").Value = "  -2.83%  "

$ws.Range("B44").Value = "CoreDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D44").Value = "'3.44"
$ws.Range("E44").Value = "  +8.20%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "'2.72"
$ws.Range("E45").Value = "  +2.16%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").Value = "'0.346"
$ws.Range("E46").Value = "  +0.70%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0429"
$ws.Range("E47").Value = "  +3.62%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "'3.28"
$ws.Range("E48").Value = "  +5.20%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.130"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.62"
$ws.Range("E50").Value = "  -0.31%  "

$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").Value = "'0.999"
$ws.Range("E51").Value = "  -0.18%  "

